# Update "想去人数" (want-to-go count) figures in the 苏州-漫展信息 workbook.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1134
$ws1.Range("F6").Value  = 4
$ws1.Range("F8").Value  = 251
$ws1.Range("F9").Value  = 394
$ws1.Range("F10").Value = 1023
$ws1.Range("F11").Value = 15
$ws1.Range("F12").Value = 520
$ws1.Range("F13").Value = 545
$ws1.Range("F15").Value = 12898
$ws1.Range("F17").Value = 5299
$ws1.Range("F18").Value = 5537

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 138

# --- Sheet "全部类型" (all types, combined list) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1134
$ws4.Range("F6").Value  = 4
$ws4.Range("F8").Value  = 251
$ws4.Range("F9").Value  = 394
$ws4.Range("F10").Value = 1023
$ws4.Range("F11").Value = 15
$ws4.Range("F12").Value = 520
$ws4.Range("F13").Value = 545
$ws4.Range("F15").Value = 12898
$ws4.Range("F16").Value = 138
$ws4.Range("F19").Value = 5299
$ws4.Range("F20").Value = 5537
